$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$ws.Range("K2").Value = "2025-11-07T03:02:39.709539+00:00"
$ws.Range("K3").Value = "2025-11-07T03:02:39.709576+00:00"
$ws.Range("K4").Value = "2025-11-07T03:02:39.709596+00:00"
$ws.Range("K5").Value = "2025-11-07T03:02:42.082762+00:00"
$ws.Range("K6").Value = "2025-11-07T03:02:42.082795+00:00"
$ws.Range("K7").Value = "2025-11-07T03:02:42.082816+00:00"
$ws.Range("K8").Value = "2025-11-07T03:02:44.267908+00:00"
$ws.Range("K9").Value = "2025-11-07T03:02:46.638705+00:00"
$ws.Range("K10").Value = "2025-11-07T03:02:46.638735+00:00"
$ws.Range("K11").Value = "2025-11-07T03:02:46.638754+00:00"
$ws.Range("K12").Value = "2025-11-07T03:02:49.466407+00:00"
$ws.Range("K13").Value = "2025-11-07T03:02:49.466438+00:00"
$ws.Range("K14").Value = "2025-11-07T03:02:49.466456+00:00"
$ws.Range("K15").Value = "2025-11-07T03:02:49.466473+00:00"
$ws.Range("K16").Value = "2025-11-07T03:02:54.535477+00:00"
$ws.Range("K17").Value = "2025-11-07T03:02:56.890756+00:00"
$ws.Range("K18").Value = "2025-11-07T03:02:59.191953+00:00"
$ws.Range("K19").Value = "2025-11-07T03:02:59.191988+00:00"
$ws.Range("K20").Value = "2025-11-07T03:03:01.546961+00:00"
$ws.Range("K21").Value = "2025-11-07T03:03:04.292237+00:00"
$ws.Range("K22").Value = "2025-11-07T03:03:04.292270+00:00"
$ws.Range("K23").Value = "2025-11-07T03:03:04.292289+00:00"
$ws.Range("K24").Value = "2025-11-07T03:03:06.746167+00:00"
$ws.Range("K25").Value = "2025-11-07T03:03:06.746195+00:00"
$ws.Range("K26").Value = "2025-11-07T03:03:06.746212+00:00"
$ws.Range("K27").Value = "2025-11-07T03:03:09.118326+00:00"
$ws.Range("K28").Value = "2025-11-07T03:03:09.118356+00:00"
$ws.Range("K29").Value = "2025-11-07T03:03:09.118373+00:00"
$ws.Range("K30").Value = "2025-11-07T03:03:09.118388+00:00"
$ws.Range("K31").Value = "2025-11-07T03:03:09.118402+00:00"
$ws.Range("K32").Value = "2025-11-07T03:03:11.904464+00:00"
$ws.Range("K33").Value = "2025-11-07T03:03:14.292718+00:00"
$ws.Range("K34").Value = "2025-11-07T03:03:14.292749+00:00"
$ws.Range("K35").Value = "2025-11-07T03:03:14.292769+00:00"
$ws.Range("K36").Value = "2025-11-07T03:03:16.723060+00:00"
$ws.Range("K37").Value = "2025-11-07T03:03:16.723093+00:00"
$ws.Range("K38").Value = "2025-11-07T03:03:16.723112+00:00"
$ws.Range("K39").Value = "2025-11-07T03:03:19.109452+00:00"
$ws.Range("K40").Value = "2025-11-07T03:03:19.109483+00:00"
$ws.Range("K41").Value = "2025-11-07T03:03:19.109501+00:00"
$ws.Range("K42").Value = "2025-11-07T03:03:19.109517+00:00"
$ws.Range("K43").Value = "2025-11-07T03:03:19.109534+00:00"
$ws.Range("K44").Value = "2025-11-07T03:03:19.109549+00:00"
$ws.Range("K45").Value = "2025-11-07T03:03:19.109566+00:00"
$ws.Range("K46").Value = "2025-11-07T03:03:19.109581+00:00"
$ws.Range("K47").Value = "2025-11-07T03:03:21.502555+00:00"
$ws.Range("K48").Value = "2025-11-07T03:03:21.502586+00:00"
$ws.Range("K49").Value = "2025-11-07T03:03:27.128353+00:00"
$ws.Range("K50").Value = "2025-11-07T03:03:27.128382+00:00"
$ws.Range("K51").Value = "2025-11-07T03:03:37.760233+00:00"
$ws.Range("K52").Value = "2025-11-07T03:03:37.760265+00:00"
